# AFDP-9108 Fix Consultation Module Admin Issues
# Fix change consultation status workflow row:
#  - Column D (Start a Workflow Process?) should hold real Boolean values
#    instead of text "true"/"false".
#  - Row 24 (Change Consultation Status) should drive the
#    "change_consultation_status" process, not "change_case_status".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 17 is the "Default Workflow" row -> should start with no workflow (false)
$ws.Range("D17").Value = $false

# Rows 18-24 all start a workflow -> should be boolean true
$ws.Range("D18").Value = $true
$ws.Range("D19").Value = $true
$ws.Range("D20").Value = $true
$ws.Range("D21").Value = $true
$ws.Range("D22").Value = $true
$ws.Range("D23").Value = $true
$ws.Range("D24").Value = $true

# Row 24 "Change Consultation Status" rule should reference the
# change_consultation_status file type/process, not change_case_status.
$ws.Range("C24").Value = "change_consultation_status"
